# BOM adjusted to numbers agreed with BB, ready to order
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Quantities ("Stück") changed:
#   CO2-Sensor / Temp-Feuchtigkeitssensor / Lora Node all shared the same
#   "5" quantity -> now "7"
$ws.Range("B3").Value = 7
$ws.Range("B4").Value = 7
$ws.Range("B5").Value = 7

#   Lora Antennenkabel quantity "1" -> "2"
$ws.Range("B6").Value = 2

# Re-stamp the "Summe" total row (C7:D7) formatting, matching the look
# the rest of the sheet already uses.
$ws.Range("C7:D7").Style = "Normal"

# Leave the cursor where the author left it after finishing the edit.
$ws.Range("A8").Select()
